$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": rows for 98f3bff0 (row 4) and d548ec63 (row 5) move from
# "Ready for handoff" to "Handed back: in sync with en-US" for both the
# zh-cn and de-de columns (E, F). The generate-date column (G) keeps its
# existing text.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value2 = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value2 = "Handed back: in sync with en-US"
$wsOverview.Range("E5").Value2 = "Handed back: in sync with en-US"
$wsOverview.Range("F5").Value2 = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": rows 4 (98f3bff0...) and 5 (d548ec63...) get fully handed
# back - Status changes, Latest Target File + Latest Handback File get
# populated, and Latest Handback DateTime is stamped.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C4").Value2 = "Handed back: in sync with en-US"
$wsZh.Range("J4").Value2 = "98f3bff0-7c12-40dd-b8eb-8bd22155eae4.8e1072d00cfe2c8c370927a86e0b73fa1971abca.zh-cn.xlf"
$wsZh.Range("K4").Value2 = "2016-08-28 12:25:47"

$wsZh.Range("C5").Value2 = "Handed back: in sync with en-US"
$wsZh.Range("J5").Value2 = "d548ec63-03d8-4c14-a3a4-6df3a5762f4c.1848d0db4de25a49668ed589f517947c757569e6.zh-cn.xlf"
$wsZh.Range("K5").Value2 = "2016-08-28 12:25:47"

$wsZh.Range("I4").Value2 = "98f3bff0-7c12-40dd-b8eb-8bd22155eae4.md"
$wsZh.Range("I5").Value2 = "d548ec63-03d8-4c14-a3a4-6df3a5762f4c.md"

# Rebuild the hyperlinks collection in row order so the new "Latest Target
# File" links for rows 4/5 land right after their row's existing hyperlink,
# matching how Excel lays out newly-added hyperlinks alongside the rest.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20902bde2b74d45631a5c6772b2533a6aabf985e/e2e/6078e45d-6b8d-4f2a-b087-85fa4331d5cc.md", "", "", "6078e45d-6b8d-4f2a-b087-85fa4331d5cc.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b5938527bfac04c3aff735a669c60b18a12e0b1c/e2e/6078e45d-6b8d-4f2a-b087-85fa4331d5cc.md", "", "", "6078e45d-6b8d-4f2a-b087-85fa4331d5cc.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20902bde2b74d45631a5c6772b2533a6aabf985e/e2e/7e314904-3a85-4fe1-a353-d521f0df730f.md", "", "", "7e314904-3a85-4fe1-a353-d521f0df730f.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b5938527bfac04c3aff735a669c60b18a12e0b1c/e2e/7e314904-3a85-4fe1-a353-d521f0df730f.md", "", "", "7e314904-3a85-4fe1-a353-d521f0df730f.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afcb1225c7aefe1b91b0437468fbd632a4a62f9a/e2e/98f3bff0-7c12-40dd-b8eb-8bd22155eae4.md", "", "", "98f3bff0-7c12-40dd-b8eb-8bd22155eae4.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b5938527bfac04c3aff735a669c60b18a12e0b1c/e2e/98f3bff0-7c12-40dd-b8eb-8bd22155eae4.md", "", "", "98f3bff0-7c12-40dd-b8eb-8bd22155eae4.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afcb1225c7aefe1b91b0437468fbd632a4a62f9a/e2e/d548ec63-03d8-4c14-a3a4-6df3a5762f4c.md", "", "", "d548ec63-03d8-4c14-a3a4-6df3a5762f4c.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b5938527bfac04c3aff735a669c60b18a12e0b1c/e2e/d548ec63-03d8-4c14-a3a4-6df3a5762f4c.md", "", "", "d548ec63-03d8-4c14-a3a4-6df3a5762f4c.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": same shape of change as zh-cn, different handback datetime
# and target-repo URLs.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C4").Value2 = "Handed back: in sync with en-US"
$wsDe.Range("J4").Value2 = "98f3bff0-7c12-40dd-b8eb-8bd22155eae4.8e1072d00cfe2c8c370927a86e0b73fa1971abca.de-de.xlf"
$wsDe.Range("K4").Value2 = "2016-08-28 12:25:53"

$wsDe.Range("C5").Value2 = "Handed back: in sync with en-US"
$wsDe.Range("J5").Value2 = "d548ec63-03d8-4c14-a3a4-6df3a5762f4c.1848d0db4de25a49668ed589f517947c757569e6.de-de.xlf"
$wsDe.Range("K5").Value2 = "2016-08-28 12:25:53"

$wsDe.Range("I4").Value2 = "98f3bff0-7c12-40dd-b8eb-8bd22155eae4.md"
$wsDe.Range("I5").Value2 = "d548ec63-03d8-4c14-a3a4-6df3a5762f4c.md"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20902bde2b74d45631a5c6772b2533a6aabf985e/e2e/6078e45d-6b8d-4f2a-b087-85fa4331d5cc.md", "", "", "6078e45d-6b8d-4f2a-b087-85fa4331d5cc.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7c35c5c46923c0780a061696a9a99ce867a1e991/e2e/6078e45d-6b8d-4f2a-b087-85fa4331d5cc.md", "", "", "6078e45d-6b8d-4f2a-b087-85fa4331d5cc.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20902bde2b74d45631a5c6772b2533a6aabf985e/e2e/7e314904-3a85-4fe1-a353-d521f0df730f.md", "", "", "7e314904-3a85-4fe1-a353-d521f0df730f.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7c35c5c46923c0780a061696a9a99ce867a1e991/e2e/7e314904-3a85-4fe1-a353-d521f0df730f.md", "", "", "7e314904-3a85-4fe1-a353-d521f0df730f.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afcb1225c7aefe1b91b0437468fbd632a4a62f9a/e2e/98f3bff0-7c12-40dd-b8eb-8bd22155eae4.md", "", "", "98f3bff0-7c12-40dd-b8eb-8bd22155eae4.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7c35c5c46923c0780a061696a9a99ce867a1e991/e2e/98f3bff0-7c12-40dd-b8eb-8bd22155eae4.md", "", "", "98f3bff0-7c12-40dd-b8eb-8bd22155eae4.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afcb1225c7aefe1b91b0437468fbd632a4a62f9a/e2e/d548ec63-03d8-4c14-a3a4-6df3a5762f4c.md", "", "", "d548ec63-03d8-4c14-a3a4-6df3a5762f4c.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7c35c5c46923c0780a061696a9a99ce867a1e991/e2e/d548ec63-03d8-4c14-a3a4-6df3a5762f4c.md", "", "", "d548ec63-03d8-4c14-a3a4-6df3a5762f4c.md") | Out-Null

Write-Output "Handback report generated"
